$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A ("Id") values: replace the numeric ids with generated
# random-id strings (poiji now reads this sheet to build the payload),
# except row 3 which becomes a different numeric id.
$ws.Range("A2").Value = "RandomId_6"
$ws.Range("A3").Value = 32453
$ws.Range("A4").Value = "RandomId_5"
$ws.Range("A5").Value = "RandomId_7"
$ws.Range("A6").Value = "RandomId"
